# Add a "Date" column as the new column A, pushing the existing survey
# columns one position to the right, and populate it with the submission
# date/time for each response row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before the current column A. Excel shifts every
# other column (headers + data) one letter to the right automatically.
$ws.Columns("A:A").Insert()

# Row 2 on this sheet carries the second-level header ("Response",
# "Open-Ended Response", ...). Label the new column "Date" there.
$ws.Range("A2").Value = "Date"

# Submission date/time for each data row (3-23), as Excel serial
# date-time numbers (days since 1899-12-30).
$dates = @(
    43843,
    43819,
    43817,
    43796,
    43795,
    43788,
    43776.503946759258,
    43775.374432870369,
    43762.409861111111,
    43762.407175925924,
    43760,
    43749.232719907406,
    43748,
    43740.354398148149,
    43739.537731481483,
    43725.755115740743,
    43720.744293981479,
    43719.629212962966,
    43718.555925925924,
    43718.542303240742,
    43706.742766203701
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 3
    $cell = $ws.Range("A$row")
    $cell.Value = $dates[$i]
    $cell.NumberFormat = "m/d/yy h:mm"
}

# Match the widened "Date" column seen in the saved workbook.
$ws.Columns("A:A").ColumnWidth = 11.5

# Restore the cursor to where the editor left it.
$ws.Range("D12").Select() | Out-Null
